$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill row 11, columns C through M with the same value/style as B11 ("x")
$ws.Range("C11:M11").Value = "x"

# Update the active selection to C19 (as reflected in the saved view state)
$ws.Range("C19").Select()
